$wb = $excel.ActiveWorkbook

# Rename the "Preschool " sheet (trailing space) to "Preschool"
$ws = $wb.Worksheets.Item("Preschool ")
$ws.Name = "Preschool"

# Make "Preschool" the active sheet/tab (it was "PrePreschool" before)
$ws.Activate()

# Set the selected cell on the now-active "Preschool" sheet
$ws.Range("O15").Select()
